# Applies scheduled-runner value updates to Sheets/Golem_Profits.xlsx
# Generated from the canonical OOXML diff: per-row H..N (current market averages,
# leve sale prices, and computed profit) refreshed with new pricing snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 85.333336
$ws.Range("I38").Value = 82.59999999999999
$ws.Range("J38").Value = 99
$ws.Range("K38").Value = 247.8
$ws.Range("L38").Value = 297
$ws.Range("M38").Value = 124.2
$ws.Range("N38").Value = -1041
# Row 103
$ws.Range("H103").Value = 999.5
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 999.5
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2998.5
$ws.Range("M103").Value = $null
$ws.Range("N103").Value = -4170.5
# Row 131
$ws.Range("H131").Value = 3692.6667
$ws.Range("I131").Value = 1054.5
$ws.Range("J131").Value = 8969
$ws.Range("K131").Value = 3163.5
$ws.Range("L131").Value = 26907
$ws.Range("M131").Value = 1876.5
$ws.Range("N131").Value = -36987
# Row 137
$ws.Range("H137").Value = 1058
$ws.Range("I137").Value = 1069.6
$ws.Range("K137").Value = 3208.8
$ws.Range("M137").Value = -658.7999999999997

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 3707
$ws.Range("I97").Value = 2510
$ws.Range("J97").Value = 4305.5
$ws.Range("K97").Value = 2510
$ws.Range("L97").Value = 4305.5
$ws.Range("M97").Value = -2014
$ws.Range("N97").Value = -5297.5
# Row 102
$ws.Range("H102").Value = 1473.6666
$ws.Range("I102").Value = 1473.6666
$ws.Range("K102").Value = 1473.6666
$ws.Range("M102").Value = 148.3334
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = $null
$ws.Range("N110").Value = $null

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 43471
$ws.Range("I26").Value = 43471
$ws.Range("K26").Value = 43471
$ws.Range("M26").Value = -43179
# Row 33
$ws.Range("H33").Value = 4673.6665
$ws.Range("J33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("N33").Value = -5672
# Row 107
$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 1100
$ws.Range("K107").Value = 1100
$ws.Range("M107").Value = 820

$ws = $wb.Worksheets.Item("CRP")
# Row 51
$ws.Range("H51").Value = 34444.332
$ws.Range("J51").Value = 29166.5
$ws.Range("L51").Value = 29166.5
$ws.Range("N51").Value = -30638.5
# Row 60
$ws.Range("H60").Value = 27942.25
$ws.Range("J60").Value = 32892.332
$ws.Range("L60").Value = 32892.332
$ws.Range("N60").Value = -33914.332
# Row 61
$ws.Range("H61").Value = 34444.332
$ws.Range("J61").Value = 29166.5
$ws.Range("L61").Value = 29166.5
$ws.Range("N61").Value = -29862.5
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
# Row 74
$ws.Range("H74").Value = 82178.5
$ws.Range("J74").Value = 82178.5
$ws.Range("L74").Value = 82178.5
$ws.Range("N74").Value = -83926.5
# Row 77
$ws.Range("H77").Value = 82178.5
$ws.Range("J77").Value = 82178.5
$ws.Range("L77").Value = 246535.5
$ws.Range("N77").Value = -255271.5
# Row 107
$ws.Range("H107").Value = 196.42857
$ws.Range("I107").Value = 212.66667
$ws.Range("K107").Value = 212.66667
$ws.Range("M107").Value = 1707.33333
# Row 122
$ws.Range("H122").Value = 2012
$ws.Range("I122").Value = 2012
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6036
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3586
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 5974.2
$ws.Range("J3").Value = 4968.25
$ws.Range("L3").Value = 14904.75
$ws.Range("N3").Value = -15128.75
# Row 68
$ws.Range("H68").Value = 2199.2
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 2499
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 7497
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -9119
# Row 71
$ws.Range("H71").Value = 2199.2
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 2499
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 22491
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -30603
# Row 134
$ws.Range("H134").Value = 967.25
$ws.Range("I134").Value = 967.25
$ws.Range("K134").Value = 2901.75
$ws.Range("M134").Value = 2168.25

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 1401
$ws.Range("I107").Value = 1400
$ws.Range("J107").Value = 1403
$ws.Range("K107").Value = 1400
$ws.Range("L107").Value = 1403
$ws.Range("M107").Value = 520
$ws.Range("N107").Value = -5243

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1960.2
$ws.Range("I16").Value = 1960.2
$ws.Range("K16").Value = 1960.2
$ws.Range("M16").Value = -1790.2
# Row 39
$ws.Range("H39").Value = 15414.75
$ws.Range("I39").Value = 553
$ws.Range("J39").Value = 60000
$ws.Range("K39").Value = 553
$ws.Range("L39").Value = 60000
$ws.Range("M39").Value = -93
$ws.Range("N39").Value = -60920
# Row 40
$ws.Range("H40").Value = 25282.818
$ws.Range("I40").Value = 21077.125
$ws.Range("J40").Value = 36498
$ws.Range("K40").Value = 21077.125
$ws.Range("L40").Value = 36498
$ws.Range("M40").Value = -20941.125
$ws.Range("N40").Value = -36770
# Row 55
$ws.Range("H55").Value = 523
$ws.Range("I55").Value = 397.75
$ws.Range("J55").Value = 648.25
$ws.Range("K55").Value = 397.75
$ws.Range("L55").Value = 648.25
$ws.Range("M55").Value = -224.75
$ws.Range("N55").Value = -994.25
# Row 68
$ws.Range("H68").Value = 3501
$ws.Range("I68").Value = 3501
$ws.Range("K68").Value = 3501
$ws.Range("M68").Value = -2752
# Row 71
$ws.Range("H71").Value = 3501
$ws.Range("I71").Value = 3501
$ws.Range("K71").Value = 17505
$ws.Range("M71").Value = -13761

$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null
# Row 137
$ws.Range("H137").Value = 95000
$ws.Range("J137").Value = 95000
$ws.Range("L137").Value = 95000
$ws.Range("N137").Value = -105200
